$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new submission ("Олег Бажайкин") needs to be inserted as row 2 (the data
# rows are stored with most-recent submission first), pushing every
# existing data row down by one. We do this by copying values upward-to-
# downward shift (processing bottom-to-top so sources aren't clobbered
# before they're read), then writing the new row's data into row 2.

for ($r = 14; $r -ge 2; $r--) {
    $dest = $r + 1
    $ws.Cells.Item($dest, 1).Value = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($dest, 2).Value = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($dest, 3).Value = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($dest, 4).Value = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($dest, 5).Value = $ws.Cells.Item($r, 5).Value2
    $ws.Cells.Item($dest, 6).Value = $ws.Cells.Item($r, 6).Value2
    $ws.Cells.Item($dest, 7).Value = $ws.Cells.Item($r, 7).Value2
}

$ws.Cells.Item(2, 1).Value = "Олег"
$ws.Cells.Item(2, 2).Value = "Бажайкин"
$ws.Cells.Item(2, 3).Value = 79521656455
$ws.Cells.Item(2, 4).Value = "bashay.oleg@mail.ru"
$ws.Cells.Item(2, 5).Value = 21
$ws.Cells.Item(2, 6).Value = "Студенты"
$ws.Cells.Item(2, 7).Value = "01.06.2023 11:31"
